$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper scratch cell used to force text (non-numeric) values through
# PasteSpecial(values-only) so numeric-looking strings like "298.30" stay
# text instead of being auto-converted to numbers by COM Value assignment.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

$ws.Range("D2").Value = "22.933.44"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.571.33"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  +0.36%  "
$helper.Value = "1.003"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.27%  "
$helper.Value = "298.30"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -1.48%  "
$helper.Value = "0.3730"
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -1.01%  "
$helper.Value = "0.3541"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -2.94%  "
$helper.Value = "49.79"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +1.75%  "
$helper.Value = "1.003"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +0.18%  "
$helper.Value = "1.204"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -5.23%  "
$helper.Value = "0.07910"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -2.13%  "
$helper.Value = "21.61"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -6.24%  "
$helper.Value = "6.383"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -3.01%  "
$helper.Value = "7.220"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -4.61%  "
$helper.Value = "0.00001210"
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -4.45%  "
$ws.Range("D17").Value = "1.577.36"
$ws.Range("E17").Value = "  -2.04%  "
$helper.Value = "91.31"
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.28%  "
$helper.Value = "0.06722"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -0.81%  "
$helper.Value = "17.58"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -4.15%  "
$helper.Value = "1.003"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +0.27%  "
$helper.Value = "6.320"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -3.90%  "
$ws.Range("D23").Value = "22.930.42"
$ws.Range("E23").Value = "  -1.18%  "
$helper.Value = "12.52"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -4.27%  "
$helper.Value = "2.368"
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +0.55%  "
$helper.Value = "2.789"
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -4.31%  "
$helper.Value = "20.43"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -3.14%  "
$helper.Value = "146.79"
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -2.41%  "
$helper.Value = "5.167"
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -1.48%  "
$helper.Value = "130.53"
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -1.45%  "
$helper.Value = "2.336"
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -3.32%  "
$helper.Value = "6.453"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -7.14%  "
$ws.Range("D33").Value = "1.755.72"
$ws.Range("E33").Value = "  -1.81%  "
$helper.Value = "0.9208"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -6.03%  "
$helper.Value = "0.07264"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -6.07%  "
$helper.Value = "0.02642"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -5.10%  "
$helper.Value = "0.08703"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$helper.Value = "9.816"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$helper.Value = "0.2450"
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -3.98%  "
$helper.Value = "5.919"
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -5.84%  "
$helper.Value = "1.333"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -4.68%  "
$helper.Value = "0.6785"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -5.07%  "
$helper.Value = "11.63"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -8.89%  "
$helper.Value = "14.61"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -8.21%  "
$ws.Range("E45").Value = "  +0.18%  "
$helper.Value = "0.6267"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -5.13%  "
$helper.Value = "3.953"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -0.74%  "
$helper.Value = "2.220"
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -3.52%  "
$helper.Value = "129.92"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -0.97%  "
$helper.Value = "0.07827"
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -2.39%  "
$helper.Value = "1.177"
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.74%  "

$helper.Clear()
$excel.CutCopyMode = $false
